$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 376.6
$ws.Range("I9").Value = 395.66666
$ws.Range("J9").Value = 348
$ws.Range("K9").Value = 395.66666
$ws.Range("L9").Value = 348
$ws.Range("M9").Value = -226.66666
$ws.Range("N9").Value = -686

$ws.Range("H17").Value = 611.5238000000001
$ws.Range("J17").Value = 611.5238000000001
$ws.Range("L17").Value = 1834.5714
$ws.Range("N17").Value = -2170.5714

$ws.Range("H115").Value = 458.42856
$ws.Range("I115").Value = 458.42856
$ws.Range("K115").Value = 1375.28568
$ws.Range("M115").Value = 191.71432

$ws.Range("H116").Value = 35924120
$ws.Range("I116").Value = 45143812
$ws.Range("J116").Value = 29416106
$ws.Range("K116").Value = 45143812
$ws.Range("L116").Value = 29416106
$ws.Range("M116").Value = -45140370
$ws.Range("N116").Value = -29422990

$ws.Range("H121").Value = 3477.739
$ws.Range("J121").Value = 3581.5
$ws.Range("L121").Value = 10744.5
$ws.Range("N121").Value = -14238.5

$ws.Range("H132").Value = 121104.56
$ws.Range("I132").Value = 365571.56
$ws.Range("K132").Value = 1096714.68
$ws.Range("M132").Value = -1094184.68

$ws.Range("H137").Value = 3802.487
$ws.Range("I137").Value = 3492
$ws.Range("J137").Value = 3924.4644
$ws.Range("K137").Value = 10476
$ws.Range("L137").Value = 11773.3932
$ws.Range("M137").Value = -7926
$ws.Range("N137").Value = -16873.3932

$ws.Range("H138").Value = 5768.3267
$ws.Range("I138").Value = 2243.9333
$ws.Range("J138").Value = 7323.206
$ws.Range("K138").Value = 6731.7999
$ws.Range("L138").Value = 21969.618
$ws.Range("M138").Value = -1591.7999
$ws.Range("N138").Value = -32249.618

$ws.Range("H141").Value = 5488.6665
$ws.Range("I141").Value = 5698
$ws.Range("J141").Value = 3500
$ws.Range("K141").Value = 17094
$ws.Range("L141").Value = 10500
$ws.Range("M141").Value = -11914
$ws.Range("N141").Value = -20860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 791977.3
$ws.Range("I2").Value = 999388
$ws.Range("J2").Value = 3816.8
$ws.Range("K2").Value = 999388
$ws.Range("L2").Value = 3816.8
$ws.Range("M2").Value = -999275
$ws.Range("N2").Value = -4042.8

$ws.Range("H32").Value = 3264.6104
$ws.Range("I32").Value = 1781.5735
$ws.Range("K32").Value = 1781.5735
$ws.Range("M32").Value = -1494.5735

$ws.Range("H45").Value = 3100
$ws.Range("I45").Value = 2502.75
$ws.Range("J45").Value = 3697.25
$ws.Range("K45").Value = 2502.75
$ws.Range("L45").Value = 3697.25
$ws.Range("M45").Value = -2125.75
$ws.Range("N45").Value = -4451.25

$ws.Range("H74").Value = 5246
$ws.Range("I74").Value = 1830.7778
$ws.Range("J74").Value = 6643.136
$ws.Range("K74").Value = 1830.7778
$ws.Range("L74").Value = 6643.136
$ws.Range("M74").Value = -956.7778000000001
$ws.Range("N74").Value = -8391.136

$ws.Range("H77").Value = 5246
$ws.Range("I77").Value = 1830.7778
$ws.Range("J77").Value = 6643.136
$ws.Range("K77").Value = 9153.889000000001
$ws.Range("L77").Value = 33215.68
$ws.Range("M77").Value = -4785.889000000001
$ws.Range("N77").Value = -41951.68

$ws.Range("H97").Value = 740.08
$ws.Range("I97").Value = 653.1905
$ws.Range("J97").Value = 1196.25
$ws.Range("K97").Value = 653.1905
$ws.Range("L97").Value = 1196.25
$ws.Range("M97").Value = -157.1905
$ws.Range("N97").Value = -2188.25

$ws.Range("H110").Value = 932790.5
$ws.Range("I110").Value = 1204876
$ws.Range("J110").Value = 7700
$ws.Range("K110").Value = 1204876
$ws.Range("L110").Value = 7700
$ws.Range("M110").Value = -1202831
$ws.Range("N110").Value = -11790

$ws.Range("H116").Value = 791977.3
$ws.Range("I116").Value = 999388
$ws.Range("J116").Value = 3816.8
$ws.Range("K116").Value = 999388
$ws.Range("L116").Value = 3816.8
$ws.Range("M116").Value = -997094
$ws.Range("N116").Value = -8404.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 791977.3
$ws.Range("I3").Value = 999388
$ws.Range("J3").Value = 3816.8
$ws.Range("K3").Value = 999388
$ws.Range("L3").Value = 3816.8
$ws.Range("M3").Value = -999274
$ws.Range("N3").Value = -4044.8

$ws.Range("H99").Value = 1390675.2
$ws.Range("I99").Value = 1489866.2
$ws.Range("K99").Value = 1489866.2
$ws.Range("M99").Value = -1488368.2

$ws.Range("H105").Value = 37039260
$ws.Range("I105").Value = 90910744
$ws.Range("J105").Value = 2615.875
$ws.Range("K105").Value = 90910744
$ws.Range("L105").Value = 2615.875
$ws.Range("M105").Value = -90908997
$ws.Range("N105").Value = -6109.875

$ws.Range("H107").Value = 1137.2727
$ws.Range("I107").Value = 1163.875
$ws.Range("J107").Value = 1066.3334
$ws.Range("K107").Value = 1163.875
$ws.Range("L107").Value = 1066.3334
$ws.Range("M107").Value = 756.125
$ws.Range("N107").Value = -4906.3334

$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2295.29
$ws.Range("I31").Value = 2351.2036
$ws.Range("J31").Value = 2229.652
$ws.Range("K31").Value = 2351.2036
$ws.Range("L31").Value = 2229.652
$ws.Range("M31").Value = -2056.2036
$ws.Range("N31").Value = -2819.652

$ws.Range("H34").Value = 2295.29
$ws.Range("I34").Value = 2351.2036
$ws.Range("J34").Value = 2229.652
$ws.Range("K34").Value = 2351.2036
$ws.Range("L34").Value = 2229.652
$ws.Range("M34").Value = -2149.2036
$ws.Range("N34").Value = -2633.652

$ws.Range("H107").Value = 759361.1
$ws.Range("I107").Value = 1136825.2
$ws.Range("K107").Value = 1136825.2
$ws.Range("M107").Value = -1134905.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 264.2857
$ws.Range("I22").Value = 264.2857
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 792.8571000000001
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -623.8571000000001
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 264.2857
$ws.Range("I27").Value = 264.2857
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 792.8571000000001
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -690.8571000000001
$ws.Range("N27").ClearContents()

$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 90000
$ws.Range("N106").Value = -91892

$ws.Range("H107").Value = 1979.2
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H128").Value = 347665
$ws.Range("I128").Value = 347665
$ws.Range("K128").Value = 1042995
$ws.Range("M128").Value = -1038015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 669.7143
$ws.Range("I97").Value = 467.09525
$ws.Range("J97").Value = 1277.5714
$ws.Range("K97").Value = 467.09525
$ws.Range("L97").Value = 1277.5714
$ws.Range("M97").Value = 28.90474999999998
$ws.Range("N97").Value = -2269.5714

$ws.Range("H107").Value = 2801661
$ws.Range("I107").Value = 4329510.5
$ws.Range("J107").Value = 603.3333
$ws.Range("K107").Value = 4329510.5
$ws.Range("L107").Value = 603.3333
$ws.Range("M107").Value = -4327590.5
$ws.Range("N107").Value = -4443.3333

$ws.Range("H122").Value = 28396068
$ws.Range("I122").Value = 1578573
$ws.Range("J122").Value = 45461748
$ws.Range("K122").Value = 4735719
$ws.Range("L122").Value = 136385244
$ws.Range("M122").Value = -4733269
$ws.Range("N122").Value = -136390144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4264.13
$ws.Range("I132").Value = 4264.13
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12792.39
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10262.39
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 98000
$ws.Range("J135").Value = 98000
$ws.Range("L135").Value = 98000
$ws.Range("N135").Value = -108140

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10026828
$ws.Range("I132").Value = 28976.195
$ws.Range("J132").Value = 125002130
$ws.Range("K132").Value = 86928.58499999999
$ws.Range("L132").Value = 375006390
$ws.Range("M132").Value = -84398.58499999999
$ws.Range("N132").Value = -375011450

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("M139").ClearContents()
